$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 142857650
$ws.Range("I12").Value = 700
$ws.Range("J12").Value = 250000350
$ws.Range("K12").Value = 700
$ws.Range("L12").Value = 250000350
$ws.Range("M12").Value = -530
$ws.Range("N12").Value = -250000690

# Row 18
$ws.Range("H18").Value = 642
$ws.Range("I18").Value = 615.6667
$ws.Range("K18").Value = 615.6667
$ws.Range("M18").Value = -331.6667

# Row 69
$ws.Range("H69").Value = 3802.175
$ws.Range("I69").Value = 3753.0588
$ws.Range("J69").Value = 3838.4783
$ws.Range("K69").Value = 11259.1764
$ws.Range("L69").Value = 11515.4349
$ws.Range("M69").Value = -10385.1764
$ws.Range("N69").Value = -13263.4349

# Row 72
$ws.Range("H72").Value = 3802.175
$ws.Range("I72").Value = 3753.0588
$ws.Range("J72").Value = 3838.4783
$ws.Range("K72").Value = 33777.5292
$ws.Range("L72").Value = 34546.3047
$ws.Range("M72").Value = -29409.5292
$ws.Range("N72").Value = -43282.3047

# Row 129
$ws.Range("H129").Value = 4543.478
$ws.Range("I129").Value = 800
$ws.Range("J129").Value = 4713.636
$ws.Range("K129").Value = 2400
$ws.Range("L129").Value = 14140.908
$ws.Range("M129").Value = 2600
$ws.Range("N129").Value = -24140.908

# Row 138
$ws.Range("H138").Value = 3117.724
$ws.Range("I138").Value = 2497.8823
$ws.Range("J138").Value = 3995.8333
$ws.Range("K138").Value = 7493.646900000001
$ws.Range("L138").Value = 11987.4999
$ws.Range("M138").Value = -2353.646900000001
$ws.Range("N138").Value = -22267.4999

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2072.963
$ws.Range("I32").Value = 1889.1389
$ws.Range("J32").Value = 3543.5557
$ws.Range("K32").Value = 1889.1389
$ws.Range("L32").Value = 3543.5557
$ws.Range("M32").Value = -1602.1389
$ws.Range("N32").Value = -4117.5557

# Row 102
$ws.Range("H102").Value = 1253.3334
$ws.Range("I102").Value = 1253.3334
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1253.3334
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 368.6666
$ws.Range("N102").ClearContents()

# Row 132
$ws.Range("H132").Value = 10827.27
$ws.Range("I132").Value = 10630.5
$ws.Range("J132").Value = 11056.833
$ws.Range("K132").Value = 31891.5
$ws.Range("L132").Value = 33170.499
$ws.Range("M132").Value = -29361.5
$ws.Range("N132").Value = -38230.499

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 1898.3846
$ws.Range("I134").Value = 1632.25
$ws.Range("J134").Value = 2324.2
$ws.Range("K134").Value = 4896.75
$ws.Range("L134").Value = 6972.599999999999
$ws.Range("M134").Value = -2361.75
$ws.Range("N134").Value = -12042.6

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 5411.9375
$ws.Range("I62").Value = 8225
$ws.Range("J62").Value = 2598.875
$ws.Range("K62").Value = 8225
$ws.Range("L62").Value = 2598.875
$ws.Range("M62").Value = -7601
$ws.Range("N62").Value = -3846.875

# Row 65
$ws.Range("H65").Value = 5411.9375
$ws.Range("I65").Value = 8225
$ws.Range("J65").Value = 2598.875
$ws.Range("K65").Value = 41125
$ws.Range("L65").Value = 12994.375
$ws.Range("M65").Value = -38005
$ws.Range("N65").Value = -19234.375

# Row 86
$ws.Range("H86").Value = 20001844
$ws.Range("I86").Value = 29413404
$ws.Range("J86").Value = 2277.25
$ws.Range("K86").Value = 29413404
$ws.Range("L86").Value = 2277.25
$ws.Range("M86").Value = -29412281
$ws.Range("N86").Value = -4523.25

# Row 89
$ws.Range("H89").Value = 20001844
$ws.Range("I89").Value = 29413404
$ws.Range("J89").Value = 2277.25
$ws.Range("K89").Value = 147067020
$ws.Range("L89").Value = 11386.25
$ws.Range("M89").Value = -147061404
$ws.Range("N89").Value = -22618.25

# Row 94
$ws.Range("H94").Value = 923
$ws.Range("J94").Value = 931
$ws.Range("L94").Value = 931
$ws.Range("N94").Value = -1833

$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 29.7
$ws.Range("J12").Value = 35
$ws.Range("L12").Value = 105
$ws.Range("N12").Value = -451

# Row 68
$ws.Range("H68").Value = 949.95
$ws.Range("I68").Value = 659.7
$ws.Range("J68").Value = 1143.45
$ws.Range("K68").Value = 1979.1
$ws.Range("L68").Value = 3430.35
$ws.Range("M68").Value = -1168.1
$ws.Range("N68").Value = -5052.35

# Row 71
$ws.Range("H71").Value = 949.95
$ws.Range("I71").Value = 659.7
$ws.Range("J71").Value = 1143.45
$ws.Range("K71").Value = 5937.3
$ws.Range("L71").Value = 10291.05
$ws.Range("M71").Value = -1881.3
$ws.Range("N71").Value = -18403.05

# Row 125
$ws.Range("H125").Value = 2921.3684
$ws.Range("I125").Value = 530
$ws.Range("J125").Value = 3202.7058
$ws.Range("K125").Value = 1590
$ws.Range("L125").Value = 9608.117400000001
$ws.Range("M125").Value = 3330
$ws.Range("N125").Value = -19448.1174

# Row 131
$ws.Range("H131").Value = 211.55385
$ws.Range("I131").Value = 175.85484
$ws.Range("J131").Value = 949.3333
$ws.Range("K131").Value = 527.56452
$ws.Range("L131").Value = 2847.9999
$ws.Range("M131").Value = 4512.43548
$ws.Range("N131").Value = -12927.9999

$ws = $wb.Worksheets.Item("GSM")
# Row 13
$ws.Range("H13").Value = 462.5
$ws.Range("I13").Value = 283.33334
$ws.Range("K13").Value = 283.33334
$ws.Range("M13").Value = -144.33334

# Row 107
$ws.Range("H107").Value = 208.66667
$ws.Range("I107").Value = 198
$ws.Range("J107").Value = 222
$ws.Range("K107").Value = 198
$ws.Range("L107").Value = 222
$ws.Range("M107").Value = 1722
$ws.Range("N107").Value = -4062

# Row 132
$ws.Range("H132").Value = 7547.125
$ws.Range("I132").Value = 10174.429
$ws.Range("J132").Value = 3868.9
$ws.Range("K132").Value = 30523.287
$ws.Range("L132").Value = 11606.7
$ws.Range("M132").Value = -27993.287
$ws.Range("N132").Value = -16666.7

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 2683.25
$ws.Range("I62").Value = 2863.25
$ws.Range("J62").Value = 2323.25
$ws.Range("K62").Value = 2863.25
$ws.Range("L62").Value = 2323.25
$ws.Range("M62").Value = -2239.25
$ws.Range("N62").Value = -3571.25

# Row 65
$ws.Range("H65").Value = 2683.25
$ws.Range("I65").Value = 2863.25
$ws.Range("J65").Value = 2323.25
$ws.Range("K65").Value = 14316.25
$ws.Range("L65").Value = 11616.25
$ws.Range("M65").Value = -11196.25
$ws.Range("N65").Value = -17856.25

# Row 96
$ws.Range("H96").Value = 1410
$ws.Range("I96").Value = 1736
$ws.Range("J96").Value = 1084
$ws.Range("K96").Value = 1736
$ws.Range("L96").Value = 1084
$ws.Range("M96").Value = -363
$ws.Range("N96").Value = -3830

# Row 107
$ws.Range("H107").Value = 343.14285
$ws.Range("I107").Value = 240.4
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 721.2
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 1198.8
$ws.Range("N107").Value = -5640

# Row 132
$ws.Range("H132").Value = 156266.4
$ws.Range("I132").Value = 183719.17
$ws.Range("J132").Value = 3315.2856
$ws.Range("K132").Value = 551157.51
$ws.Range("L132").Value = 9945.856800000001
$ws.Range("M132").Value = -548627.51
$ws.Range("N132").Value = -15005.8568
